$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.229934581839814
$ws.Range("C2").Value = 0.4056676705255313
$ws.Range("D2").Value = 0.01145294848190481
$ws.Range("E2").Value = 0.05052024881184558
$ws.Range("F2").Value = 3.628100717984182
$ws.Range("I2").Value = 2.057488557387408
$ws.Range("J2").Value = 0.1112702309760154
$ws.Range("L2").Value = 0.4919442578562752
$ws.Range("N2").Value = 2.150983187293836
$ws.Range("B3").Value = 3.117137602726928
$ws.Range("C3").Value = 0.37525515478481
$ws.Range("D3").Value = 0.01110465608398847
$ws.Range("E3").Value = 0.05058848726883303
$ws.Range("F3").Value = 3.6154679144136
$ws.Range("I3").Value = 2.05879922382119
$ws.Range("J3").Value = 0.1119254443326474
$ws.Range("L3").Value = 0.4853832162082625
$ws.Range("N3").Value = 2.174632267203641
$ws.Range("B4").Value = 3.049638617032031
$ws.Range("C4").Value = 0.3568102768785479
$ws.Range("D4").Value = 0.01088656952492428
$ws.Range("E4").Value = 0.0506343506248752
$ws.Range("F4").Value = 3.609672594614764
$ws.Range("I4").Value = 2.060670984191795
$ws.Range("J4").Value = 0.1123519967548861
$ws.Range("L4").Value = 0.4815950824093278
$ws.Range("N4").Value = 2.189888701198647
$ws.Range("B5").Value = 3.022573837985476
$ws.Range("C5").Value = 0.3493507527680038
$ws.Range("D5").Value = 0.01079659825481194
$ws.Range("E5").Value = 0.05065403867714025
$ws.Range("F5").Value = 3.607802938907724
$ws.Range("I5").Value = 2.06170141970027
$ws.Range("J5").Value = 0.1125319266524238
$ws.Range("L5").Value = 0.4801118100636614
$ws.Range("N5").Value = 2.196290577318745
$ws.Range("B6").Value = 3.018106405460458
$ws.Range("C6").Value = 0.34811552594158
$ws.Range("D6").Value = 0.0107815911958884
$ws.Range("E6").Value = 0.05065736820472744
$ws.Range("F6").Value = 3.607522164122074
$ws.Range("I6").Value = 2.061888671596108
$ws.Range("J6").Value = 0.1125621729406845
$ws.Range("L6").Value = 0.4798691626974261
$ws.Range("N6").Value = 2.197364754984271
$ws.Range("B7").Value = 3.04927182478616
$ws.Range("C7").Value = 0.3567094454076027
$ws.Range("D7").Value = 0.01088536063688217
$ws.Range("E7").Value = 0.05063461210027898
$ws.Range("F7").Value = 3.60964538936021
$ws.Range("I7").Value = 2.060683798037203
$ws.Range("J7").Value = 0.1123543986160369
$ws.Range("L7").Value = 0.4815748338534291
$ws.Range("N7").Value = 2.189974291445893
$ws.Range("B8").Value = 3.190676688735437
$ws.Range("C8").Value = 0.3951336593515009
$ws.Range("D8").Value = 0.01133371497453872
$ws.Range("E8").Value = 0.05054295539599907
$ws.Range("F8").Value = 3.623337128173162
$ws.Range("I8").Value = 2.057718632064613
$ws.Range("J8").Value = 0.1114911210117953
$ws.Range("L8").Value = 0.4896321047581154
$ws.Range("N8").Value = 2.158984514163457
$ws.Range("B9").Value = 3.48197650906593
$ws.Range("C9").Value = 0.4723260261688438
$ws.Range("D9").Value = 0.01218093158030165
$ws.Range("E9").Value = 0.05039461646815602
$ws.Range("F9").Value = 3.665808928163017
$ws.Range("I9").Value = 2.060401108645536
$ws.Range("J9").Value = 0.1099902342024475
$ws.Range("L9").Value = 0.5073421517624155
$ws.Range("N9").Value = 2.104063065235742
$ws.Range("B10").Value = 3.704630519717171
$ws.Range("C10").Value = 0.530209519198479
$ws.Range("D10").Value = 0.01278610244690093
$ws.Range("E10").Value = 0.05030470144844124
$ws.Range("F10").Value = 3.706630235277146
$ws.Range("I10").Value = 2.067598579541325
$ws.Range("J10").Value = 0.1090040058790036
$ws.Range("L10").Value = 0.5215237486772537
$ws.Range("N10").Value = 2.067292283782322
$ws.Range("B11").Value = 3.807820368240073
$ws.Range("C11").Value = 0.5568069704143568
$ws.Range("D11").Value = 0.01305814622149093
$ws.Range("E11").Value = 0.05026792373418676
$ws.Range("F11").Value = 3.727310351989189
$ws.Range("I11").Value = 2.072018281737925
$ws.Range("J11").Value = 0.108580514932866
$ws.Range("L11").Value = 0.5282308357392509
$ws.Range("N11").Value = 2.051344766703167
$ws.Range("B12").Value = 3.847170789188738
$ws.Range("C12").Value = 0.5669177241861121
$ws.Range("D12").Value = 0.01316073416429653
$ws.Range("E12").Value = 0.05025458913176695
$ws.Range("F12").Value = 3.735446427469128
$ws.Range("I12").Value = 2.073857450026694
$ws.Range("J12").Value = 0.1084237577117957
$ws.Range("L12").Value = 0.530807500847672
$ws.Range("N12").Value = 2.045418295349123
$ws.Range("B13").Value = 3.838683740470856
$ws.Range("C13").Value = 0.5647384532504134
$ws.Range("D13").Value = 0.01313865848432982
$ws.Range("E13").Value = 0.0502574346464475
$ws.Range("F13").Value = 3.733680590397711
$ws.Range("I13").Value = 2.073453977717705
$ws.Range("J13").Value = 0.1084573577537657
$ws.Range("L13").Value = 0.5302509307156242
$ws.Range("N13").Value = 2.046689656928201
$ws.Range("B14").Value = 3.811052240728941
$ws.Range("C14").Value = 0.557638005416436
$ws.Range("D14").Value = 0.01306659457197767
$ws.Range("E14").Value = 0.05026681482209083
$ws.Range("F14").Value = 3.727973591158417
$ws.Range("I14").Value = 2.07216626928421
$ws.Range("J14").Value = 0.1085675461185875
$ws.Range("L14").Value = 0.5284420809518906
$ws.Range("N14").Value = 2.050854935273982
$ws.Range("B15").Value = 3.794162947887514
$ws.Range("C15").Value = 0.5532938578671178
$ws.Range("D15").Value = 0.01302239860807219
$ws.Range("E15").Value = 0.05027263756380029
$ws.Range("F15").Value = 3.724517649825572
$ws.Range("I15").Value = 2.071399090528374
$ws.Range("J15").Value = 0.1086355095512328
$ws.Range("L15").Value = 0.5273389070943466
$ws.Range("N15").Value = 2.053420951952681
$ws.Range("B16").Value = 3.697925202751264
$ws.Range("C16").Value = 0.5284767187317811
$ws.Range("D16").Value = 0.01276826159234368
$ws.Range("E16").Value = 0.05030718789053434
$ws.Range("F16").Value = 3.705321333699089
$ws.Range("I16").Value = 2.067332860482907
$ws.Range("J16").Value = 0.1090321874901754
$ws.Range("L16").Value = 0.5210905766147107
$ws.Range("N16").Value = 2.068350213968877
$ws.Range("B17").Value = 3.639374640396227
$ws.Range("C17").Value = 0.513320772219231
$ws.Range("D17").Value = 0.01261155118439383
$ws.Range("E17").Value = 0.05032943929370137
$ws.Range("F17").Value = 3.694086576233985
$ws.Range("I17").Value = 2.065132321629747
$ws.Range("J17").Value = 0.1092819729946388
$ws.Range("L17").Value = 0.5173229823737842
$ws.Range("N17").Value = 2.077708804767319
$ws.Range("B18").Value = 3.605877004758895
$ws.Range("C18").Value = 0.5046284676063806
$ws.Range("D18").Value = 0.01252110713171106
$ws.Range("E18").Value = 0.05034262603157957
$ws.Range("F18").Value = 3.687823180579784
$ws.Range("I18").Value = 2.063974404798245
$ws.Range("J18").Value = 0.1094280105793466
$ws.Range("L18").Value = 0.5151800362679069
$ws.Range("N18").Value = 2.083164983633687
$ws.Range("B19").Value = 3.594566012996552
$ws.Range("C19").Value = 0.5016896764706757
$ws.Range("D19").Value = 0.01249043020283658
$ws.Range("E19").Value = 0.05034715755605523
$ws.Range("F19").Value = 3.685736556365555
$ws.Range("I19").Value = 2.063600840445488
$ws.Range("J19").Value = 0.1094778632380535
$ws.Range("L19").Value = 0.5144586048294428
$ws.Range("N19").Value = 2.085024940853721
$ws.Range("B20").Value = 3.645588904854037
$ws.Range("C20").Value = 0.5149315571771353
$ws.Range("D20").Value = 0.01262826491273472
$ws.Range("E20").Value = 0.05032703041071307
$ws.Range("F20").Value = 3.695261976429805
$ws.Range("I20").Value = 2.065355412732458
$ws.Range("J20").Value = 0.1092551379167244
$ws.Range("L20").Value = 0.517721557234637
$ws.Range("N20").Value = 2.076704971025244
$ws.Range("B21").Value = 3.819160819835133
$ws.Range("C21").Value = 0.5597225188838593
$ws.Range("D21").Value = 0.01308777282789642
$ws.Range("E21").Value = 0.05026404356833747
$ws.Range("F21").Value = 3.729641586396241
$ws.Range("I21").Value = 2.0725400018042
$ws.Range("J21").Value = 0.1085350832147416
$ws.Range("L21").Value = 0.5289723837995695
$ws.Range("N21").Value = 2.04962843550571
$ws.Range("B22").Value = 3.934201533063117
$ws.Range("C22").Value = 0.5892229033023
$ws.Range("D22").Value = 0.01338560466891359
$ws.Range("E22").Value = 0.05022633016882017
$ws.Range("F22").Value = 3.753888931910978
$ws.Range("I22").Value = 2.078200720091772
$ws.Range("J22").Value = 0.1080855223217867
$ws.Range("L22").Value = 0.5365401895165149
$ws.Range("N22").Value = 2.032588318796044
$ws.Range("B23").Value = 3.872655280549907
$ws.Range("C23").Value = 0.5734570198066535
$ws.Range("D23").Value = 0.01322686060340672
$ws.Range("E23").Value = 0.05024614291515661
$ws.Range("F23").Value = 3.740784449413383
$ws.Range("I23").Value = 2.075090908019831
$ws.Range("J23").Value = 0.1083235388166823
$ws.Range("L23").Value = 0.5324814421737614
$ws.Range("N23").Value = 2.041622798559338
$ws.Range("B24").Value = 3.642778925567427
$ws.Range("C24").Value = 0.5142032557759535
$ws.Range("D24").Value = 0.01262070972625118
$ws.Range("E24").Value = 0.05032811823910222
$ws.Range("F24").Value = 3.694729968814357
$ws.Range("I24").Value = 2.065254219393822
$ws.Range("J24").Value = 0.1092672624778146
$ws.Range("L24").Value = 0.5175412896108185
$ws.Range("N24").Value = 2.07715856818826
$ws.Range("B25").Value = 3.401662852223637
$ws.Range("C25").Value = 0.4512413664301107
$ws.Range("D25").Value = 0.01195493421160165
$ws.Range("E25").Value = 0.05043139226766868
$ws.Range("F25").Value = 3.652637135319623
$ws.Range("I25").Value = 2.058760982826939
$ws.Range("J25").Value = 0.110375768888141
$ws.Range("L25").Value = 0.5023460782419562
$ws.Range("N25").Value = 2.118293131375822
